# Insert a new bullet-list paragraph right after the paragraph that reads
# "Case Classes: Verticles Event Bus Pattern Matching (Dispatcher / Signatures)."
# The new paragraph keeps the same list formatting (numId 8 / ilvl 0,
# ind left=600 hanging=360) as its neighbours and carries an explicit
# "no underline" paragraph-mark run property, per the target OOXML.

$d = $word.ActiveDocument

# Locate the anchor paragraph by its exact text.
$anchorText = "Case Classes: Verticles Event Bus Pattern Matching (Dispatcher / Signatures)."
$paras = $d.Paragraphs
$anchor = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $candidate = $paras.Item($i)
    if ($candidate.Range.Text -like "*$anchorText*") {
        $anchor = $candidate
        break
    }
}

if ($anchor -eq $null) {
    throw "Anchor paragraph not found"
}

# Create a fresh empty paragraph right after the anchor (inherits the
# anchor's pPr, i.e. the same numbering / indent).
$tail = $anchor.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# The empty paragraph we just created is now the paragraph right after
# the anchor; re-fetch it by index so we have a fresh Range.
$newIndex = $anchor.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)

# Replace its (empty) contents with the exact target OOXML, giving full
# control over pPr/rPr (including the paragraph-mark "no underline").
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:ind w:left="600" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Reactive / Event Driven: Verticles DIDs (Distributed IDs) distributed patterns routing registry. Resource / Applicable graph logs. Rx Facade.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml)
